$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 58
$ws.Range("H58").Value = 21611.295
$ws.Range("J58").Value = 25585.28
$ws.Range("L58").Value = 76755.84
$ws.Range("N58").Value = -77055.84
# Row 132
$ws.Range("H132").Value = 2942911.5
$ws.Range("I132").Value = 3031833
$ws.Range("J132").Value = 8500
$ws.Range("K132").Value = 9095499
$ws.Range("L132").Value = 25500
$ws.Range("M132").Value = -9092969
$ws.Range("N132").Value = -30560
# Row 135
$ws.Range("H135").Value = 749.4318
$ws.Range("I135").Value = 533.475
$ws.Range("K135").Value = 4801.275000000001
$ws.Range("M135").Value = -2266.275000000001
# Row 137
$ws.Range("H137").Value = 2059.8438
$ws.Range("I137").Value = 2051.6
$ws.Range("K137").Value = 6154.799999999999
$ws.Range("M137").Value = -3604.799999999999
# Row 138
$ws.Range("H138").Value = 4437.787
$ws.Range("I138").Value = 1544.1818
$ws.Range("J138").Value = 6070.077
$ws.Range("K138").Value = 4632.5454
$ws.Range("L138").Value = 18210.231
$ws.Range("M138").Value = 507.4546
$ws.Range("N138").Value = -28490.231
# Row 141
$ws.Range("H141").Value = 412321.28
$ws.Range("I141").Value = 1525.6666
$ws.Range("J141").Value = 1644708.1
$ws.Range("K141").Value = 4576.9998
$ws.Range("L141").Value = 4934124.300000001
$ws.Range("M141").Value = 603.0002000000004
$ws.Range("N141").Value = -4944484.300000001

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 8929867
$ws.Range("I2").Value = 13889750
$ws.Range("J2").Value = 2078.5
$ws.Range("K2").Value = 13889750
$ws.Range("L2").Value = 2078.5
$ws.Range("M2").Value = -13889637
$ws.Range("N2").Value = -2304.5
# Row 40
$ws.Range("H40").Value = 70031
$ws.Range("J40").Value = 70031
$ws.Range("L40").Value = 70031
$ws.Range("N40").Value = -70383
# Row 45
$ws.Range("H45").Value = 1319.6923
$ws.Range("I45").Value = 1044.6471
$ws.Range("J45").Value = 3190
$ws.Range("K45").Value = 1044.6471
$ws.Range("L45").Value = 3190
$ws.Range("M45").Value = -667.6470999999999
$ws.Range("N45").Value = -3944
# Row 116
$ws.Range("H116").Value = 8929867
$ws.Range("I116").Value = 13889750
$ws.Range("J116").Value = 2078.5
$ws.Range("K116").Value = 13889750
$ws.Range("L116").Value = 2078.5
$ws.Range("M116").Value = -13887456
$ws.Range("N116").Value = -6666.5
# Row 132
$ws.Range("H132").Value = 24392718
$ws.Range("I132").Value = 34484428
$ws.Range("K132").Value = 103453284
$ws.Range("M132").Value = -103450754

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 8929867
$ws.Range("I3").Value = 13889750
$ws.Range("J3").Value = 2078.5
$ws.Range("K3").Value = 13889750
$ws.Range("L3").Value = 2078.5
$ws.Range("M3").Value = -13889636
$ws.Range("N3").Value = -2306.5
# Row 35
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
# Row 46
$ws.Range("H46").Value = 5000
$ws.Range("I46").Value = 5000
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 5000
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -4702
$ws.Range("N46").ClearContents()
# Row 107
$ws.Range("H107").Value = 3498.5386
$ws.Range("I107").Value = 2259.1667
$ws.Range("J107").Value = 4560.857
$ws.Range("K107").Value = 2259.1667
$ws.Range("L107").Value = 4560.857
$ws.Range("M107").Value = -339.1667000000002
$ws.Range("N107").Value = -8400.857
# Row 134
$ws.Range("H134").Value = 2674.4849
$ws.Range("I134").Value = 2106.077
$ws.Range("J134").Value = 4785.7144
$ws.Range("K134").Value = 6318.231000000001
$ws.Range("L134").Value = 14357.1432
$ws.Range("M134").Value = -3783.231000000001
$ws.Range("N134").Value = -19427.1432

$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Range("H3").Value = 55334.332
$ws.Range("J3").Value = 70001.5
$ws.Range("L3").Value = 70001.5
$ws.Range("N3").Value = -70227.5
# Row 16
$ws.Range("H16").Value = 2749.5557
$ws.Range("J16").Value = 3024.1428
$ws.Range("L16").Value = 3024.1428
$ws.Range("N16").Value = -3598.1428
# Row 31
$ws.Range("H31").Value = 2699.76
$ws.Range("I31").Value = 1751.6364
$ws.Range("J31").Value = 3444.7144
$ws.Range("K31").Value = 1751.6364
$ws.Range("L31").Value = 3444.7144
$ws.Range("M31").Value = -1456.6364
$ws.Range("N31").Value = -4034.7144
# Row 34
$ws.Range("H34").Value = 2699.76
$ws.Range("I34").Value = 1751.6364
$ws.Range("J34").Value = 3444.7144
$ws.Range("K34").Value = 1751.6364
$ws.Range("L34").Value = 3444.7144
$ws.Range("M34").Value = -1549.6364
$ws.Range("N34").Value = -3848.7144
# Row 113
$ws.Range("H113").Value = 2749.5557
$ws.Range("J113").Value = 3024.1428
$ws.Range("L113").Value = 3024.1428
$ws.Range("N113").Value = -7364.1428
# Row 122
$ws.Range("H122").Value = 3538.353
$ws.Range("I122").Value = 2832.625
$ws.Range("J122").Value = 4165.6665
$ws.Range("K122").Value = 8497.875
$ws.Range("L122").Value = 12496.9995
$ws.Range("M122").Value = -6047.875
$ws.Range("N122").Value = -17396.9995
# Row 132
$ws.Range("H132").Value = 4650
$ws.Range("I132").Value = 3350
$ws.Range("J132").Value = 7250
$ws.Range("K132").Value = 10050
$ws.Range("L132").Value = 21750
$ws.Range("M132").Value = -7520
$ws.Range("N132").Value = -26810

$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 75
$ws.Range("I23").Value = 65
$ws.Range("J23").Value = 95
$ws.Range("K23").Value = 195
$ws.Range("L23").Value = 285
$ws.Range("M23").Value = 40
$ws.Range("N23").Value = -755
# Row 122
$ws.Range("H122").Value = 1434.9412
$ws.Range("J122").Value = 2574
$ws.Range("L122").Value = 23166
$ws.Range("N122").Value = -28066

$ws = $wb.Worksheets.Item("GSM")
# Row 46
$ws.Range("H46").Value = 9923
$ws.Range("J46").Value = 9923
$ws.Range("L46").Value = 9923
$ws.Range("N46").Value = -10235
# Row 102
$ws.Range("H102").Value = 114957
$ws.Range("I102").Value = 2233.3333
$ws.Range("J102").Value = 171318.83
$ws.Range("K102").Value = 2233.3333
$ws.Range("L102").Value = 171318.83
$ws.Range("M102").Value = -611.3332999999998
$ws.Range("N102").Value = -174562.83

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 76927750
$ws.Range("I61").Value = 125002860
$ws.Range("K61").Value = 125002860
$ws.Range("M61").Value = -125002658
# Row 113
$ws.Range("H113").Value = 76927750
$ws.Range("I113").Value = 125002860
$ws.Range("K113").Value = 125002860
$ws.Range("M113").Value = -125000690

$ws = $wb.Worksheets.Item("WVR")
# Row 40
$ws.Range("H40").Value = 70028
$ws.Range("J40").Value = 70028
$ws.Range("L40").Value = 70028
$ws.Range("N40").Value = -70326
# Row 132
$ws.Range("H132").Value = 8401.843999999999
$ws.Range("I132").Value = 1788.4762
$ws.Range("K132").Value = 5365.4286
$ws.Range("M132").Value = -2835.4286
